$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "department" column (C) values on rows 2-7 to reflect the
# specific department for each course instead of the generic school name.
$ws.Range("C2").Value = "Early Childhood"
$ws.Range("C3").Value = "Early Childhood"
$ws.Range("C4").Value = "Ageing Support"
$ws.Range("C5").Value = "Ageing Support"
$ws.Range("C6").Value = "Ageing Support"
$ws.Range("C7").Value = "Community Services"
